# Updated symbol list refresh: new "Price" (col D) readings and "Hora" (col G)
# bumped from 4 -> 5 for every data row (2-51). Values are written with a
# leading "'" so Excel stores them as text (quotePrefix), matching the
# original inline-string cell types instead of being reinterpreted as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'282.19"
$ws.Range("G2").Value = "'5"

$ws.Range("D3").Value = "'20.62"
$ws.Range("G3").Value = "'5"

$ws.Range("D4").Value = "'6.206"
$ws.Range("G4").Value = "'5"

$ws.Range("D5").Value = "'0.06173"
$ws.Range("G5").Value = "'5"

$ws.Range("D6").Value = "'3.587"
$ws.Range("G6").Value = "'5"

$ws.Range("G7").Value = "'5"

$ws.Range("D8").Value = "'1.497"
$ws.Range("G8").Value = "'5"

$ws.Range("D9").Value = "'0.8190"
$ws.Range("G9").Value = "'5"

$ws.Range("D10").Value = "'0.01383"
$ws.Range("G10").Value = "'5"

$ws.Range("D11").Value = "'0.1634"
$ws.Range("G11").Value = "'5"

$ws.Range("D12").Value = "'0.08392"
$ws.Range("G12").Value = "'5"

$ws.Range("D13").Value = "'0.03512"
$ws.Range("G13").Value = "'5"

$ws.Range("D14").Value = "'0.03218"
$ws.Range("G14").Value = "'5"

$ws.Range("D15").Value = "'0.09152"
$ws.Range("G15").Value = "'5"

$ws.Range("D16").Value = "'3.701"
$ws.Range("G16").Value = "'5"

$ws.Range("D17").Value = "'0.001643"
$ws.Range("G17").Value = "'5"

$ws.Range("D18").Value = "'0.04720"
$ws.Range("G18").Value = "'5"

$ws.Range("D19").Value = "'0.006487"
$ws.Range("G19").Value = "'5"

$ws.Range("D20").Value = "'0.006174"
$ws.Range("G20").Value = "'5"

$ws.Range("G21").Value = "'5"

$ws.Range("G22").Value = "'5"

$ws.Range("G23").Value = "'5"

$ws.Range("G24").Value = "'5"

$ws.Range("D25").Value = "'0.3356"
$ws.Range("G25").Value = "'5"

$ws.Range("D26").Value = "'0.1232"
$ws.Range("G26").Value = "'5"

$ws.Range("G27").Value = "'5"

$ws.Range("G28").Value = "'5"

$ws.Range("G29").Value = "'5"

$ws.Range("G30").Value = "'5"

$ws.Range("G31").Value = "'5"

$ws.Range("G32").Value = "'5"

$ws.Range("G33").Value = "'5"

$ws.Range("G34").Value = "'5"

$ws.Range("G35").Value = "'5"

$ws.Range("G36").Value = "'5"

$ws.Range("G37").Value = "'5"

$ws.Range("G38").Value = "'5"

$ws.Range("G39").Value = "'5"

$ws.Range("D40").Value = "'0.04706"
$ws.Range("G40").Value = "'5"

$ws.Range("D41").Value = "'0.007199"
$ws.Range("G41").Value = "'5"

$ws.Range("D42").Value = "'0.1101"
$ws.Range("G42").Value = "'5"

$ws.Range("D43").Value = "'0.003495"
$ws.Range("G43").Value = "'5"

$ws.Range("D44").Value = "'0.01141"
$ws.Range("G44").Value = "'5"

$ws.Range("D45").Value = "'0.00006533"
$ws.Range("G45").Value = "'5"

$ws.Range("G46").Value = "'5"

$ws.Range("D47").Value = "'1.051"
$ws.Range("G47").Value = "'5"

$ws.Range("D48").Value = "'0.002841"
$ws.Range("G48").Value = "'5"

$ws.Range("G49").Value = "'5"

$ws.Range("D50").Value = "'0.01242"
$ws.Range("G50").Value = "'5"

$ws.Range("G51").Value = "'5"
